$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows for the "house upgrade" feature (key | en | fr)
$rows = @(
    @("UPGRADE_TITLE",       "UPGRADE", "AMÉLIORER"),
    @("UPGRADE_STONE_LABEL", "Stone",   "Pierre"),
    @("UPGRADE_GOLD_LABEL",  "Gold",    "Or"),
    @("UPGRADE_BUTTON",      "UPGRADE", "AMÉLIORER")
)

$startRow = 19
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Range("F20").Select()
